$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column (D) whose new values are ambiguous numeric-looking
# strings (e.g. "0.999", "593.11"). Excel would silently convert a plain
# assignment of such text into a real number, so we briefly force a text
# number format while writing the value, then restore the default "Normal"
# cell style (no explicit number format), matching the original workbook.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D26", "D29", "D30", "D32", "D33", "D34", "D37", "D43", "D45", "D46", "D48", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated Price (D) and Volume(1h) (E) figures for this run
$ws.Range("D2").Value = '61.300.84'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '2.937.98'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '593.11'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("D6").Value = '146.08'
$ws.Range("E6").Value = '  +1.38%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("D9").Value = '6.98'
$ws.Range("E9").Value = '  +4.44%  '
$ws.Range("D10").Value = '0.144'
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").Value = '0.442'
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '0.0000227'
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("D13").Value = '33.86'
$ws.Range("E13").Value = '  +1.38%  '
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '3.419.69'
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("D16").Value = '61.117.18'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '6.74'
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").Value = '2.936.74'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").Value = '437.92'
$ws.Range("E19").Value = '  +2.54%  '
$ws.Range("D20").Value = '13.48'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").Value = '0.680'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("D23").Value = '81.76'
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("D24").Value = '11.03'
$ws.Range("E24").Value = '  +1.34%  '
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").Value = '11.89'
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("E28").Value = '  +3.95%  '
$ws.Range("D29").Value = '2.61'
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = '7.04'
$ws.Range("E30").Value = '  -2.25%  '
$ws.Range("E31").Value = '  +3.86%  '
$ws.Range("D32").Value = '26.77'
$ws.Range("E32").Value = '  +1.74%  '
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("D34").Value = '0.0₃0877'
$ws.Range("E34").Value = '  +2.23%  '
$ws.Range("E35").Value = '  +0.89%  '
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("D37").Value = '3.03'
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("E41").Value = '  +2.17%  '
$ws.Range("E42").Value = '  -2.10%  '
$ws.Range("D43").Value = '378.62'
$ws.Range("E43").Value = '  +1.13%  '
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").Value = '2.696.08'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").Value = '132.90'
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("D48").Value = '24.08'
$ws.Range("E48").Value = '  +0.26%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("D51").Value = '0.126'
$ws.Range("E51").Value = '  +1.55%  '

# Restore the default (unstyled/General) appearance on the Price cells we touched
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
